$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165461182594299
$ws.Range("B1").Value = 2.371947288513184
$ws.Range("D1").Value = 2.391100168228149
$ws.Range("E1").Value = 1.215088605880737
